# Updated cryptos list with GitHub Actions
# Applies the latest price / 1h-volume-change snapshot to the crypto table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    # Writes a value while forcing it to be stored as text, so that
    # numeric-looking strings (e.g. "1.00", "68.139.85") are not silently
    # converted into real numbers by Excel. The cell's original style is
    # preserved (saved before, restored after) so no visible formatting
    # changes are introduced.
    param($Ws, $Addr, $Val)
    $rng = $Ws.Range($Addr)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $Val
    $rng.Style = $origStyle
}

# --- Row 39 / 40: TheGraph and Bittensor swapped ranking position ---
Set-TextValue $ws "B39" "Bittensor"
$ws.Range("C39").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue $ws "D39" "456.17"
$ws.Range("E39").Value = "  +4.21%  "

Set-TextValue $ws "B40" "TheGraph"
$ws.Range("C40").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
Set-TextValue $ws "D40" "0.325"
$ws.Range("E40").Value = "  -8.50%  "

# --- Price (D) and Volume(1h) (E) updates ---
Set-TextValue $ws "D2" "68.139.85"
$ws.Range("E2").Value = "  -3.23%  "
Set-TextValue $ws "D3" "3.816.83"
$ws.Range("E3").Value = "  +1.52%  "
$ws.Range("E4").Value = "  +0.02%  "
Set-TextValue $ws "D5" "594.84"
$ws.Range("E5").Value = "  -4.03%  "
Set-TextValue $ws "D6" "173.08"
$ws.Range("E6").Value = "  -5.06%  "
Set-TextValue $ws "D7" "3.814.86"
$ws.Range("E7").Value = "  +1.49%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("E9").Value = "  -1.05%  "
Set-TextValue $ws "D10" "0.160"
$ws.Range("E10").Value = "  -4.58%  "
Set-TextValue $ws "D11" "6.29"
$ws.Range("E11").Value = "  -1.63%  "
$ws.Range("E12").Value = "  -3.07%  "
Set-TextValue $ws "D13" "38.35"
$ws.Range("E13").Value = "  -4.79%  "
Set-TextValue $ws "D14" "0.0000246"
$ws.Range("E14").Value = "  -4.21%  "
Set-TextValue $ws "D15" "4.454.19"
$ws.Range("E15").Value = "  +1.59%  "
Set-TextValue $ws "D16" "3.814.79"
$ws.Range("E16").Value = "  +1.63%  "
Set-TextValue $ws "D17" "68.197.37"
$ws.Range("E17").Value = "  -3.27%  "
$ws.Range("E18").Value = "  -4.49%  "
Set-TextValue $ws "D19" "7.18"
$ws.Range("E19").Value = "  -5.24%  "
Set-TextValue $ws "D20" "16.11"
$ws.Range("E20").Value = "  -2.91%  "
Set-TextValue $ws "D21" "490.24"
$ws.Range("E21").Value = "  -3.28%  "
Set-TextValue $ws "D22" "9.30"
$ws.Range("E22").Value = "  +0.79%  "
Set-TextValue $ws "D23" "0.738"
$ws.Range("E23").Value = "  +2.00%  "
Set-TextValue $ws "D24" "85.47"
$ws.Range("E24").Value = "  -2.23%  "
Set-TextValue $ws "D25" "2.38"
$ws.Range("E25").Value = "  -8.79%  "
Set-TextValue $ws "D26" "0.0000139"
$ws.Range("E26").Value = "  +2.55%  "
Set-TextValue $ws "D27" "12.25"
$ws.Range("E27").Value = "  -6.68%  "
Set-TextValue $ws "D28" "10.18"
$ws.Range("E28").Value = "  -10.75%  "
$ws.Range("E29").Value = "  -0.05%  "
$ws.Range("E30").Value = "  -0.02%  "
$ws.Range("E31").Value = "  -2.01%  "
Set-TextValue $ws "D32" "32.90"
$ws.Range("E32").Value = "  +7.37%  "
Set-TextValue $ws "D33" "7.72"
$ws.Range("E33").Value = "  -2.69%  "
$ws.Range("E34").Value = "  -4.16%  "
Set-TextValue $ws "D35" "0.999"
$ws.Range("E35").Value = "  +0.03%  "
$ws.Range("E36").Value = "  -4.01%  "
Set-TextValue $ws "D37" "0.137"
$ws.Range("E37").Value = "  -2.03%  "
Set-TextValue $ws "D38" "5.79"
$ws.Range("E38").Value = "  -5.42%  "
Set-TextValue $ws "D41" "49.00"
$ws.Range("E41").Value = "  -2.18%  "
$ws.Range("E42").Value = "  -4.04%  "
Set-TextValue $ws "D43" "2.90"
$ws.Range("E43").Value = "  -7.06%  "
Set-TextValue $ws "D44" "8.28"
$ws.Range("E44").Value = "  -4.11%  "
Set-TextValue $ws "D45" "41.69"
$ws.Range("E45").Value = "  -7.56%  "
Set-TextValue $ws "D46" "2.845.11"
$ws.Range("E46").Value = "  -3.89%  "
$ws.Range("E47").Value = "  +0.03%  "
$ws.Range("E48").Value = "  -3.12%  "
Set-TextValue $ws "D49" "138.57"
$ws.Range("E49").Value = "  +0.85%  "
Set-TextValue $ws "D50" "26.33"
$ws.Range("E50").Value = "  -3.64%  "
Set-TextValue $ws "D51" "23.48"
$ws.Range("E51").Value = "  +9.24%  "
